$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.53849
$ws.Range("E2").Value = 1390
$ws.Range("F2").Value = 0.21302
$ws.Range("I2").Value = 0.06018
$ws.Range("K2").Value = 2.08018
$ws.Range("M2").Value = 1.88554
$ws.Range("N2").Value = 0.00101
